$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.061.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.821.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5980'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = '  -6.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2746'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.10'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07605'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.823.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.719'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6227'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009654'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.721.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.559'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -11.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  -6.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.835'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.904'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.25%  '
$ws.Range("E26").Value = '  -4.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.89%  '
$ws.Range("B28").Value = 'Hedera'
$ws.Range("C28").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06366'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -10.44%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.415'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.435'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.14%  '
$ws.Range("E31").Value = '  -5.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.742'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.720'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.084'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6428'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.537'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.734'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01748'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.536'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.142.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8770'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.971.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000114'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.601'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.416'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05499'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4535'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.404'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.60%  '
